$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add 4 new "select other" columns (J,K,L,M) mirroring the existing
# --- "Another select field" / "Type your answer" pair currently in H/I ---

# Copy formatting from column I (the "other" free-text column) onto the
# new columns for every existing row so borders/fills/number formats match.
$ws.Range("I1").Copy()
$ws.Range("J1:M1").PasteSpecial(-4122)
$ws.Range("I2").Copy()
$ws.Range("J2:M2").PasteSpecial(-4122)
$ws.Range("I3").Copy()
$ws.Range("J3:M3").PasteSpecial(-4122)
$ws.Range("I4:I10").Copy()
$ws.Range("J4:M10").PasteSpecial(-4122)

# Header row (row 1): two more "select field" / "other answer" pairs
$ws.Range("J1").Value = "Third select field"
$ws.Range("K1").Value = "Type your answer"
$ws.Range("L1").Value = "Fourth select field"
$ws.Range("M1").Value = "Type your answer"

# Row 2 sample data
$ws.Range("J2").Value = "Other"
$ws.Range("K2").Value = "Answer three"
$ws.Range("L2").Value = "Other"
$ws.Range("M2").Value = "Answer Four"

# Row 3 sample data (K3 stays blank, just like G3)
$ws.Range("J3").Value = "Other"
$ws.Range("L3").Value = "Other"
$ws.Range("M3").Value = "Answer Four"

# The data rows (and their blank cells) now get a solid white fill
# instead of no fill, so the borders read correctly even over gridlines.
$ws.Range("A2:M10").Interior.ColorIndex = 2

# --- Minor column width tweaks (F, G, H) and new columns I..M share I's width ---
$ws.Columns.Item(6).ColumnWidth = 22.714285714285715
$ws.Columns.Item(7).ColumnWidth = 20.42857142857143
$ws.Columns.Item(8).ColumnWidth = 26.285714285714285
$ws.Range("I1:M1").EntireColumn.ColumnWidth = 34.42857142857143
